$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for the "New2PE" column (ports on the cortex for the lift motors
# after the drive motors were moved to the brain / cortex).
$ws.Range("G1").Value = "New2PE"

# Fill in the new port values for each motor row.
$ws.Range("G2").Value = "-"
$ws.Range("G3").Value = "-"
$ws.Range("G4").Value = "-"
$ws.Range("G5").Value = "D"
$ws.Range("G6").Value = "C"
$ws.Range("G7").Value = "B"
$ws.Range("G8").Value = "A"
$ws.Range("G9").Value = "-"
$ws.Range("G10").Value = "-"
$ws.Range("G11").Value = "-"

# Match the saved selection from the authored workbook.
$ws.Range("G9").Select()
